# Weekly update: a new daily price record is inserted for
# "Terminal Hortofrutícola Agro Chillán - Repollo" at row 436. All the
# existing records from row 436 downward shift down by one row, and the
# dataset grows from A1:R473 to A1:R474.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing row 436 (and everything below it) down by one row.
$ws.Rows.Item(436).Insert()

# Populate the newly-opened row 436 with the new record.
$ws.Cells.Item(436, 1).Value  = 7
$ws.Cells.Item(436, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(436, 3).Value  = "Ñuble"
$ws.Cells.Item(436, 4).Value  = 45132
$ws.Cells.Item(436, 5).Value  = 16
$ws.Cells.Item(436, 6).Value  = 100112006
$ws.Cells.Item(436, 7).Value  = "Repollo"
$ws.Cells.Item(436, 8).Value  = "Crespo record"
$ws.Cells.Item(436, 9).Value  = "Primera"
$ws.Cells.Item(436, 10).Value = 200
$ws.Cells.Item(436, 11).Value = 1000
$ws.Cells.Item(436, 12).Value = 1000
$ws.Cells.Item(436, 13).Value = 1000
$ws.Cells.Item(436, 14).Value = "$/unidad"
$ws.Cells.Item(436, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(436, 16).Value = 1000
$ws.Cells.Item(436, 17).Value = 1
$ws.Cells.Item(436, 18).Value = "Hortaliza"
